$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.424.06"
$ws.Range("E2").Value = "  -3.90%  "
$ws.Range("D3").Value = "1.949.76"
$ws.Range("E3").Value = "  -2.72%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "321.16"
$ws.Range("E5").Value = "  -2.64%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4761"
$ws.Range("E7").Value = "  -5.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4051"
$ws.Range("E8").Value = "  -4.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "53.49"
$ws.Range("E9").Value = "  -0.97%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08478"
$ws.Range("E10").Value = "  -6.34%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.057"
$ws.Range("E11").Value = "  -5.69%  "
$ws.Range("E12").Value = "  -5.53%  "
$ws.Range("D13").Value = "1.971.25"
$ws.Range("E13").Value = "  -2.76%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.590"
$ws.Range("E14").Value = "  -5.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.173"
$ws.Range("E15").Value = "  -4.76%  "
$ws.Range("E16").Value = "  +0.33%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001075"
$ws.Range("E17").Value = "  -3.75%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "89.10"
$ws.Range("E18").Value = "  -5.73%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06598"
$ws.Range("E19").Value = "  -1.25%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.69"
$ws.Range("E20").Value = "  -5.11%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.012"
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.807"
$ws.Range("E22").Value = "  -2.86%  "
$ws.Range("D23").Value = "28.450.96"
$ws.Range("E23").Value = "  -3.86%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.50"
$ws.Range("E24").Value = "  -4.23%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.292"
$ws.Range("E25").Value = "  -0.33%  "
$ws.Range("D26").Value = "2.185.32"
$ws.Range("E26").Value = "  -2.60%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "154.42"
$ws.Range("E27").Value = "  -2.84%  "
$ws.Range("E28").Value = "  -2.87%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.955"
$ws.Range("E29").Value = "  -6.32%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.152"
$ws.Range("E30").Value = "  -6.82%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "123.72"
$ws.Range("E31").Value = "  -3.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9884"
$ws.Range("E32").Value = "  -6.79%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09562"
$ws.Range("E33").Value = "  -4.17%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.443"
$ws.Range("E34").Value = "  -8.20%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.599"
$ws.Range("E35").Value = "  -4.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.662"
$ws.Range("E36").Value = "  -3.52%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02334"
$ws.Range("E37").Value = "  -5.56%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.06226"
$ws.Range("E38").Value = "  -2.77%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.786"
$ws.Range("E39").Value = "  -5.80%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.256"
$ws.Range("E40").Value = "  -3.85%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6224"
$ws.Range("E41").Value = "  -5.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.10"
$ws.Range("E42").Value = "  -5.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.011"
$ws.Range("E43").Value = "  +0.00%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1922"
$ws.Range("E44").Value = "  -6.44%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.332"
$ws.Range("E45").Value = "  +2.12%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5964"
$ws.Range("E46").Value = "  -6.40%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.99"
$ws.Range("E47").Value = "  -3.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.056"
$ws.Range("E48").Value = "  -6.66%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.391"
$ws.Range("E49").Value = "  -3.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00000000330"
$ws.Range("E50").Value = "  -0.05%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06817"
$ws.Range("E51").Value = "  -2.52%  "
